$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7 values
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = 1

# Clear the empty inline string cells in row 7 (Q7, R7, T7)
$ws.Range("Q7").Value = ""
$ws.Range("R7").Value = ""
$ws.Range("T7").Value = ""

# Add new row 8 data
$ws.Range("A8").Value = 2
$ws.Range("B8").Value = "'7/1/2020"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 14
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.1
$ws.Range("H8").Value = 1.1
$ws.Range("I8").Value = 1.1
$ws.Range("J8").Value = 1.1
$ws.Range("K8").Value = 1.1
$ws.Range("L8").Value = 1.1
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 1.1
$ws.Range("O8").Value = 1.1
$ws.Range("P8").Value = 1.1
$ws.Range("S8").Value = "RP-7/1/2020"
$ws.Range("U8").Value = "hello"

# Update selection
$ws.Range("C11").Select()
